# "Lab 3 Complete Project 1 Server Editions" update to the CIS 232 Server Log.
# Fills in rows 11-15 of the log with the newly performed maintenance steps
# and nudges the saved view/selection to where the user left off (B16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# All five new log entries happened on the same day as the serial date
# value 42759 (2017-01-24), matching the rest of the "Date" column.
$logDate = Get-Date -Year 2017 -Month 1 -Day 24 -Hour 0 -Minute 0 -Second 0

# Rows 5-10 already carry the thin-border/date-number-format style (s="7")
# on column A. Copy that formatting onto A11:A16 instead of typing a raw
# NumberFormat string, so we reuse the existing style record rather than
# minting a new one.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A11:A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 11; Action = "Remove NIC Team 1";                                   Down = "No";  How = "N/A" },
    @{ Row = 12; Action = "Set Eth1 Static IP 10.130.26.131";                    Down = "No";  How = "N/A" },
    @{ Row = 13; Action = "Install Roles (AD DS, DHCP, DNS, IIS, Print)";        Down = "No";  How = "N/A" },
    @{ Row = 14; Action = "Remove Print Roles";                                  Down = "No";  How = "N/A" },
    @{ Row = 15; Action = "Enable roles";                                        Down = "Yes"; How = "Seconds" }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $logDate
    $ws.Range("B$r").Value = $entry.Action
    $ws.Range("C$r").Value = $entry.Down
    $ws.Range("D$r").Value = $entry.How
    $ws.Range("E$r").Value = "Admin"
    $ws.Range("F$r").Value = "Evan"
    $ws.Range("G$r").Value = "ES"
}

# Row 16 stays blank - only its date-formatted style (already applied above
# via the PasteSpecial) changes; the rest of the row keeps its old style.

# Match the saved view: scrolled so row 4 is at the top, cell B16 selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B16").Select() | Out-Null
